# Generate Report for Handback
# Update the generated timestamps in the handback-status workbook.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for aaca600f-246e-424c-8ed9-f78834fd48cd.md
$wsOverview.Range("G2").Value = "2016-08-15 18:59:58"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for aaca600f-246e-424c-8ed9-f78834fd48cd.md
$wsZhCn.Range("H2").Value = "2016-08-15 18:59:54"
$wsZhCn.Range("K2").Value = "2016-08-15 19:00:35"

# de-de sheet: "Correspond Handoff Datetime" (shares text with Overview!G2) and
# "Correspond Handback DateTime" for aaca600f-246e-424c-8ed9-f78834fd48cd.md
$wsDeDe.Range("H2").Value = "2016-08-15 18:59:58"
$wsDeDe.Range("K2").Value = "2016-08-15 19:00:42"
